$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 1.73
$ws.Range("G3").Value = 1.82
$ws.Range("H3").Value = 5.5
$ws.Range("J3").Value = 3.55
$ws.Range("K3").Value = 4
$ws.Range("P3").Value = 1.77
$ws.Range("Q3").Value = 2.1
